$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "97-53=",
    "38+37=",
    "78-61=",
    "73-33=",
    "8-2=",
    "44+3=",
    "30+26=",
    "32+40=",
    "55+11=",
    "91-68=",
    "32+44=",
    "14+44=",
    "27-9=",
    "79-73=",
    "99-32=",
    "50+2=",
    "31+43=",
    "61-57=",
    "17+71=",
    "64-50=",
    "89-47=",
    "25+71=",
    "72+13=",
    "89-7=",
    "18+33=",
    "80-36=",
    "70-65=",
    "6+65=",
    "87-28=",
    "59-13=",
    "77-59=",
    "28-16=",
    "17-16=",
    "53+41=",
    "60+34=",
    "88-59=",
    "98-64=",
    "91-64=",
    "39-32=",
    "26+47=",
    "34-28=",
    "31+18=",
    "4+77=",
    "90-64=",
    "44+48=",
    "83+13=",
    "65+7=",
    "7+58=",
    "97-82=",
    "46-11=",
    "47+38=",
    "12+1=",
    "7+43=",
    "28+9=",
    "43+29=",
    "60+28=",
    "11+32=",
    "34+48=",
    "22+69=",
    "64+33=",
    "28-8=",
    "4+91=",
    "66-15=",
    "11+74=",
    "88-11=",
    "49+35=",
    "95-37=",
    "93-15=",
    "2+78=",
    "74-62=",
    "40+53=",
    "91-65=",
    "31+35=",
    "38-19=",
    "61+32=",
    "66-24=",
    "20-0=",
    "12+63=",
    "29-6=",
    "4+95=",
    "70+16=",
    "40-5=",
    "81+12=",
    "91-9=",
    "41+18=",
    "79-18=",
    "41-33=",
    "31-16=",
    "60-23=",
    "8+40=",
    "54-41=",
    "62-51=",
    "12+13=",
    "10+31=",
    "55+5=",
    "21+74=",
    "94-16=",
    "34-18=",
    "40-22=",
    "32+28="
)

$idx = 0
foreach ($row in $t.Rows) {
    foreach ($cell in $row.Cells) {
        $r = $cell.Range
        $r.MoveEnd(1, -1) | Out-Null
        $r.Text = $values[$idx]
        $idx = $idx + 1
    }
}
Write-Output "updated $idx cells"
